# Add a new timelog entry row (row 27), mirroring the formatting of the
# previous entry row (row 26), then update the active cell / scroll state
# to reflect the new bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 26's formatting (fonts/borders/number format) down into row 27
# so the new entry matches the rest of the log.
$ws.Range("A26:B26").Copy()
$ws.Range("A27:B27").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# New timelog entry
$ws.Range("A27").Value = "3/11, 3 hrs"
$ws.Range("B27").Value = "Work on govt data in shiny, documenting, cleaning up graphs"

# Match the wrapped-text row height used for similar entries
$ws.Rows.Item(27).RowHeight = 41.4

# Leave the selection on the next empty row, like the author did after
# finishing their entry
$ws.Range("B28").Select()
